# Rimesso check a inizio_lavorazione anziche' fine_lavorazione sulla release date
# per le ricerche locali.
#
# Effetto sui dati: la riga 30 (commessa 251651 / BIMEC 4) deve essere
# riposizionata subito dopo la riga 11, e le righe 12..29 devono scalare
# di una posizione in basso (diventando le righe 13..30).
#
# Leggiamo prima tutti i valori delle righe coinvolte (12..30) in memoria,
# cosi' da poterli riscrivere nel nuovo ordine senza che le scritture
# sovrascrivano dati che devono ancora essere letti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 12
$lastRow = 30
$lastCol = "S"

# Conserva in memoria il contenuto originale di ogni riga (12..30)
$originalRows = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $originalRows[$r] = $ws.Range("A$r`:$lastCol$r").Value()
}

# Nuovo ordinamento:
#   riga 12            <- vecchia riga 30
#   riga r (13..30)     <- vecchia riga (r-1)
$ws.Range("A$firstRow`:$lastCol$firstRow").Value = $originalRows[$lastRow]

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Range("A$r`:$lastCol$r").Value = $originalRows[$r - 1]
}
